$wb = $excel.ActiveWorkbook

$urban = $wb.Worksheets.Item("mapping_urban")
$rural = $wb.Worksheets.Item("mapping_rural")

$rural.Range("B2").Value = "49.5% MUR/LWAL+CDN/H:1`n49.5% MCF/LWAL+CDL/H:1`n1% W/LWAL+CDL/H:1"
$urban.Range("B2").Value = "42.5% MUR/LWAL+CDN/H:1`n56.5% MCF/LWAL+CDL/H:1`n1% W/LWAL+CDL/H:1"
